$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '30.486.95'
$ws.Range("E2").Value = '  -0.84%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.893.34'
$ws.Range("E3").Value = '  -0.70%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.08%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '

# Row 6 - USDC
$ws.Range("E6").Value = '  +0.06%  '

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4836'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.04%  '

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2898'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.16%  '

# Row 9 - Dogecoin
$ws.Range("E9").Value = '  -1.85%  '

# Row 10 - WrappedEther
$ws.Range("D10").Value = '1.891.10'
$ws.Range("E10").Value = '  -0.86%  '

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07400'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.68%  '

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.202'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.63%  '

# Row 14 - Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.06%  '

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6627'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.18%  '

# Row 16 - WrappedBTC
$ws.Range("D16").Value = '30.462.27'
$ws.Range("E16").Value = '  -0.78%  '

# Row 17 - Avalanche
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.68%  '

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007777'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.66%  '

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '

# Row 20 - WrappedliquidstakedEther2.0
$ws.Range("D20").Value = '2.142.56'
$ws.Range("E20").Value = '  -0.60%  '

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.405'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.20%  '

# Row 22 - BinanceUSD
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.33%  '

# Row 24 - Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.210'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.81%  '

# Row 25 - Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.413'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.10%  '

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.42%  '

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.69%  '

# Row 28 - LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.943'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.06%  '

# Row 29 - Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.444'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.13%  '

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.348'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.33%  '

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09190'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.54%  '

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.058'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.46%  '

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05083'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.42%  '

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7568'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.84%  '

# Row 35 - ARBITRUM
$ws.Range("E35").Value = '  +4.62%  '

# Row 36 - HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.708'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.86%  '

# Row 37 - VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01887'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.34%  '

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.652'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.17%  '

# Row 39 - RenderToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.104'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.18%  '

# Row 40 - TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9196'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.23%  '

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.026'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.22%  '

# Row 42 - Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.35%  '

# Row 43 - TheSandbox
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4358'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.84%  '

# Row 44 - PaxDollar
$ws.Range("E44").Value = '  +0.44%  '

# Row 45 - Aptos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.651'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.52%  '

# Row 46 - NEARProtocol -> Algorand
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1339'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.48%  '

# Row 47 - Algorand -> NEARProtocol
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.598'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.60%  '

# Row 48 - Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -13.01%  '

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.959'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.59%  '

# Row 50 - Elrond
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.63%  '

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05709'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.60%  '
